# Novo grafico de rentabildiade beta
# Apply updated values to the "Negocios" (I), "Strike VS Cot." (G),
# "Real Time" (E), "TIR (%)" (F) and "Prob. Exec." (H) columns on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Cells.Item(2, 9).Value = 0
$ws.Cells.Item(3, 7).Value = -0.66
$ws.Cells.Item(3, 9).Value = 20000
$ws.Cells.Item(4, 9).Value = 0
$ws.Cells.Item(7, 7).Value = -0.61
$ws.Cells.Item(7, 9).Value = 20000
$ws.Cells.Item(8, 9).Value = 3500
$ws.Cells.Item(9, 9).Value = 81200
$ws.Cells.Item(10, 5).Value = 0.01
$ws.Cells.Item(10, 6).Value = 0.0008
$ws.Cells.Item(10, 9).Value = 65000
$ws.Cells.Item(11, 9).Value = 123200
$ws.Cells.Item(12, 7).Value = -0.52
$ws.Cells.Item(12, 9).Value = 0
$ws.Cells.Item(13, 5).Value = 0.01
$ws.Cells.Item(13, 6).Value = 0.0007
$ws.Cells.Item(13, 9).Value = 2600
$ws.Cells.Item(14, 9).Value = 11100
$ws.Cells.Item(15, 9).Value = 330900
$ws.Cells.Item(16, 5).Value = 0.02
$ws.Cells.Item(16, 6).Value = 0.0013
$ws.Cells.Item(16, 9).Value = 26400
$ws.Cells.Item(17, 5).Value = 0.02
$ws.Cells.Item(17, 6).Value = 0.0013
$ws.Cells.Item(17, 9).Value = 111400
$ws.Cells.Item(18, 5).Value = 0.02
$ws.Cells.Item(18, 6).Value = 0.0013
$ws.Cells.Item(18, 9).Value = 34800
$ws.Cells.Item(19, 5).Value = 0.02
$ws.Cells.Item(19, 6).Value = 0.0012
$ws.Cells.Item(19, 7).Value = -0.43
$ws.Cells.Item(19, 9).Value = 21400
$ws.Cells.Item(20, 8).Value = 0
$ws.Cells.Item(20, 9).Value = 0
$ws.Cells.Item(21, 9).Value = 1700
$ws.Cells.Item(22, 8).Value = 0.01
$ws.Cells.Item(22, 9).Value = 8000
$ws.Cells.Item(23, 9).Value = 14900
$ws.Cells.Item(24, 8).Value = 0.03
$ws.Cells.Item(24, 9).Value = 0
$ws.Cells.Item(25, 5).Value = 0.05
$ws.Cells.Item(25, 6).Value = 0.0029
$ws.Cells.Item(25, 8).Value = 0.04
$ws.Cells.Item(26, 8).Value = 0.08
$ws.Cells.Item(26, 9).Value = 73400
$ws.Cells.Item(27, 8).Value = 0.11
$ws.Cells.Item(27, 9).Value = 25600
$ws.Cells.Item(28, 5).Value = 0.04
$ws.Cells.Item(28, 6).Value = 0.0022
$ws.Cells.Item(28, 8).Value = 0.16
$ws.Cells.Item(29, 5).Value = 0.03
$ws.Cells.Item(29, 6).Value = 0.0016
$ws.Cells.Item(29, 8).Value = 0.21
$ws.Cells.Item(29, 9).Value = 100
$ws.Cells.Item(30, 5).Value = 0.05
$ws.Cells.Item(30, 6).Value = 0.0026
$ws.Cells.Item(30, 8).Value = 0.29
$ws.Cells.Item(30, 9).Value = 62600
$ws.Cells.Item(31, 8).Value = 0.39
$ws.Cells.Item(31, 9).Value = 500
$ws.Cells.Item(32, 5).Value = 0.06
$ws.Cells.Item(32, 6).Value = 0.0031
$ws.Cells.Item(32, 8).Value = 0.51
$ws.Cells.Item(32, 9).Value = 7000
$ws.Cells.Item(33, 8).Value = 0.67
$ws.Cells.Item(33, 9).Value = 0
$ws.Cells.Item(34, 5).Value = 0.06
$ws.Cells.Item(34, 6).Value = 0.003
$ws.Cells.Item(34, 8).Value = 0.86
$ws.Cells.Item(34, 9).Value = 443900
$ws.Cells.Item(35, 5).Value = 0.06
$ws.Cells.Item(35, 6).Value = 0.003
$ws.Cells.Item(35, 8).Value = 1.1
$ws.Cells.Item(35, 9).Value = 1000
$ws.Cells.Item(36, 8).Value = 1.39
$ws.Cells.Item(36, 9).Value = 136300
$ws.Cells.Item(37, 5).Value = 0.08
$ws.Cells.Item(37, 6).Value = 0.0039
$ws.Cells.Item(37, 8).Value = 1.74
$ws.Cells.Item(37, 9).Value = 10200
$ws.Cells.Item(38, 8).Value = 2.15
$ws.Cells.Item(38, 9).Value = 290500
$ws.Cells.Item(39, 5).Value = 0.08
$ws.Cells.Item(39, 6).Value = 0.0038
$ws.Cells.Item(39, 8).Value = 2.64
$ws.Cells.Item(39, 9).Value = 15000
$ws.Cells.Item(40, 5).Value = 0.08
$ws.Cells.Item(40, 6).Value = 0.0037
$ws.Cells.Item(40, 8).Value = 3.21
$ws.Cells.Item(40, 9).Value = 47000
$ws.Cells.Item(41, 5).Value = 0.1
$ws.Cells.Item(41, 6).Value = 0.0046
$ws.Cells.Item(41, 8).Value = 3.88
$ws.Cells.Item(41, 9).Value = 35000
$ws.Cells.Item(42, 5).Value = 0.11
$ws.Cells.Item(42, 6).Value = 0.005
$ws.Cells.Item(42, 8).Value = 4.64
$ws.Cells.Item(42, 9).Value = 360900
$ws.Cells.Item(43, 5).Value = 0.11
$ws.Cells.Item(43, 6).Value = 0.0049
$ws.Cells.Item(43, 8).Value = 5.5
$ws.Cells.Item(43, 9).Value = 8700
$ws.Cells.Item(44, 5).Value = 0.13
$ws.Cells.Item(44, 6).Value = 0.0058
$ws.Cells.Item(44, 8).Value = 6.48
$ws.Cells.Item(44, 9).Value = 209300
$ws.Cells.Item(45, 5).Value = 0.14
$ws.Cells.Item(45, 6).Value = 0.0062
$ws.Cells.Item(45, 8).Value = 7.57
$ws.Cells.Item(45, 9).Value = 66600
$ws.Cells.Item(46, 5).Value = 0.16
$ws.Cells.Item(46, 6).Value = 0.007
$ws.Cells.Item(46, 8).Value = 8.789999999999999
$ws.Cells.Item(46, 9).Value = 554600
$ws.Cells.Item(47, 5).Value = 0.14
$ws.Cells.Item(47, 6).Value = 0.006
$ws.Cells.Item(47, 7).Value = -0.17
$ws.Cells.Item(47, 8).Value = 10.13
$ws.Cells.Item(47, 9).Value = 57200
$ws.Cells.Item(48, 5).Value = 0.2
$ws.Cells.Item(48, 6).Value = 0.008500000000000001
$ws.Cells.Item(48, 8).Value = 11.59
$ws.Cells.Item(48, 9).Value = 140100
$ws.Cells.Item(49, 5).Value = 0.22
$ws.Cells.Item(49, 6).Value = 0.009299999999999999
$ws.Cells.Item(49, 8).Value = 13.18
$ws.Cells.Item(49, 9).Value = 50200
$ws.Cells.Item(50, 5).Value = 0.24
$ws.Cells.Item(50, 6).Value = 0.01
$ws.Cells.Item(50, 8).Value = 14.9
$ws.Cells.Item(50, 9).Value = 785700
$ws.Cells.Item(51, 5).Value = 0.27
$ws.Cells.Item(51, 6).Value = 0.0111
$ws.Cells.Item(51, 8).Value = 16.74
$ws.Cells.Item(51, 9).Value = 129900
$ws.Cells.Item(52, 5).Value = 0.35
$ws.Cells.Item(52, 6).Value = 0.0141
$ws.Cells.Item(52, 8).Value = 20.77
$ws.Cells.Item(52, 9).Value = 55300
$ws.Cells.Item(53, 5).Value = 0.41
$ws.Cells.Item(53, 6).Value = 0.0164
$ws.Cells.Item(53, 8).Value = 22.96
$ws.Cells.Item(53, 9).Value = 2387300
$ws.Cells.Item(54, 5).Value = 0.44
$ws.Cells.Item(54, 6).Value = 0.0174
$ws.Cells.Item(54, 8).Value = 25.24
$ws.Cells.Item(54, 9).Value = 112500
$ws.Cells.Item(55, 5).Value = 0.5
$ws.Cells.Item(55, 6).Value = 0.0196
$ws.Cells.Item(55, 8).Value = 27.61
$ws.Cells.Item(55, 9).Value = 206800
$ws.Cells.Item(56, 8).Value = 27.61
$ws.Cells.Item(57, 5).Value = 0.5600000000000001
$ws.Cells.Item(57, 6).Value = 0.0217
$ws.Cells.Item(57, 7).Value = -0.08
$ws.Cells.Item(57, 8).Value = 30.06
$ws.Cells.Item(57, 9).Value = 186100
$ws.Cells.Item(58, 5).Value = 0.62
$ws.Cells.Item(58, 6).Value = 0.0238
$ws.Cells.Item(58, 8).Value = 32.58
$ws.Cells.Item(58, 9).Value = 2005900
$ws.Cells.Item(59, 5).Value = 0.72
$ws.Cells.Item(59, 6).Value = 0.0274
$ws.Cells.Item(59, 8).Value = 35.15
$ws.Cells.Item(59, 9).Value = 348500
$ws.Cells.Item(60, 5).Value = 0.8
$ws.Cells.Item(60, 6).Value = 0.0302
$ws.Cells.Item(60, 8).Value = 37.77
$ws.Cells.Item(60, 9).Value = 46500
$ws.Cells.Item(61, 8).Value = 40.43
$ws.Cells.Item(61, 9).Value = 300200
$ws.Cells.Item(62, 5).Value = 0.95
$ws.Cells.Item(62, 6).Value = 0.0352
$ws.Cells.Item(62, 8).Value = 43.1
$ws.Cells.Item(62, 9).Value = 828300
$ws.Cells.Item(63, 5).Value = 1.2
$ws.Cells.Item(63, 6).Value = 0.0436
$ws.Cells.Item(63, 8).Value = 48.45
$ws.Cells.Item(63, 9).Value = 338200
$ws.Cells.Item(64, 8).Value = 48.45
$ws.Cells.Item(65, 5).Value = 1.33
$ws.Cells.Item(65, 6).Value = 0.0479
$ws.Cells.Item(65, 8).Value = 51.11
$ws.Cells.Item(65, 9).Value = 122300
$ws.Cells.Item(66, 5).Value = 1.47
$ws.Cells.Item(66, 6).Value = 0.0525
$ws.Cells.Item(66, 8).Value = 53.74
$ws.Cells.Item(66, 9).Value = 1215600
$ws.Cells.Item(67, 5).Value = 1.74
$ws.Cells.Item(67, 6).Value = 0.0611
$ws.Cells.Item(67, 8).Value = 58.86
$ws.Cells.Item(67, 9).Value = 84100
$ws.Cells.Item(68, 8).Value = 58.86
$ws.Cells.Item(69, 5).Value = 1.91
$ws.Cells.Item(69, 6).Value = 0.0664
$ws.Cells.Item(69, 8).Value = 61.34
$ws.Cells.Item(69, 9).Value = 39500
$ws.Cells.Item(70, 5).Value = 2.06
$ws.Cells.Item(70, 6).Value = 0.07099999999999999
$ws.Cells.Item(70, 8).Value = 63.76
$ws.Cells.Item(70, 9).Value = 582100
$ws.Cells.Item(71, 5).Value = 2.4
$ws.Cells.Item(71, 6).Value = 0.0814
$ws.Cells.Item(71, 8).Value = 68.36
$ws.Cells.Item(71, 9).Value = 10700
$ws.Cells.Item(72, 5).Value = 2.82
$ws.Cells.Item(72, 6).Value = 0.094
$ws.Cells.Item(72, 8).Value = 72.63
$ws.Cells.Item(72, 9).Value = 990500
$ws.Cells.Item(73, 5).Value = 3.21
$ws.Cells.Item(73, 6).Value = 0.1052
$ws.Cells.Item(73, 8).Value = 76.52
$ws.Cells.Item(73, 9).Value = 8000
$ws.Cells.Item(74, 8).Value = 76.52
$ws.Cells.Item(75, 5).Value = 3.36
$ws.Cells.Item(75, 6).Value = 0.1093
$ws.Cells.Item(75, 7).Value = 0.1
$ws.Cells.Item(75, 8).Value = 78.33
$ws.Cells.Item(76, 5).Value = 3.62
$ws.Cells.Item(76, 6).Value = 0.1168
$ws.Cells.Item(76, 8).Value = 80.03
$ws.Cells.Item(76, 9).Value = 32200
$ws.Cells.Item(77, 5).Value = 4.07
$ws.Cells.Item(77, 6).Value = 0.1292
$ws.Cells.Item(77, 8).Value = 83.16
$ws.Cells.Item(77, 9).Value = 11300
$ws.Cells.Item(78, 5).Value = 4.15
$ws.Cells.Item(78, 6).Value = 0.1307
$ws.Cells.Item(78, 8).Value = 84.58
$ws.Cells.Item(78, 9).Value = 200
$ws.Cells.Item(79, 5).Value = 4.55
$ws.Cells.Item(79, 6).Value = 0.1422
$ws.Cells.Item(79, 8).Value = 85.91
$ws.Cells.Item(79, 9).Value = 1100
$ws.Cells.Item(80, 8).Value = 88.31
$ws.Cells.Item(81, 5).Value = 5.36
$ws.Cells.Item(81, 6).Value = 0.1624
$ws.Cells.Item(81, 7).Value = 0.18
$ws.Cells.Item(81, 8).Value = 90.37
$ws.Cells.Item(81, 9).Value = 0
